# Refresh the cryptocurrency price (column D) and 1h change/volume (column E)
# figures for the rows whose data changed in this run, matching the source diff.
# Column D values that look numeric (e.g. "1.00", "214.91") must stay as plain
# text, exactly as stored in the original workbook, so they are written with a
# leading apostrophe to force text and then restyled back to "Normal" so no
# stray number-format style gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.680.28"
$ws.Cells.Item(2, 5).Value = "  +0.44%  "
$ws.Cells.Item(3, 4).Value = "1.642.77"
$ws.Cells.Item(3, 5).Value = "  +0.77%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'214.91"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.90%  "
$ws.Cells.Item(6, 5).Value = "  +1.28%  "
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(8, 5).Value = "  +0.74%  "
$ws.Cells.Item(9, 5).Value = "  +0.70%  "
$ws.Cells.Item(10, 5).Value = "  -0.08%  "
$ws.Cells.Item(11, 4).Value = "'0.0844"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.00%  "
$ws.Cells.Item(12, 4).Value = "1.872.67"
$ws.Cells.Item(12, 5).Value = "  +0.88%  "
$ws.Cells.Item(13, 4).Value = "1.641.99"
$ws.Cells.Item(13, 5).Value = "  +2.60%  "
$ws.Cells.Item(14, 4).Value = "'4.18"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.55%  "
$ws.Cells.Item(15, 5).Value = "  +1.44%  "
$ws.Cells.Item(16, 4).Value = "'64.95"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.56%  "
$ws.Cells.Item(17, 4).Value = "26.705.34"
$ws.Cells.Item(17, 5).Value = "  +0.05%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0744"
$ws.Cells.Item(18, 5).Value = "  +0.61%  "
$ws.Cells.Item(19, 4).Value = "'215.43"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.12%  "
$ws.Cells.Item(20, 4).Value = "'1.00"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.04%  "
$ws.Cells.Item(21, 4).Value = "'4.36"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.06%  "
$ws.Cells.Item(22, 4).Value = "'6.25"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.29%  "
$ws.Cells.Item(23, 5).Value = "  +1.41%  "
$ws.Cells.Item(24, 4).Value = "'2.24"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +14.35%  "
$ws.Cells.Item(25, 4).Value = "'145.49"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.98%  "
$ws.Cells.Item(26, 5).Value = "  +0.04%  "
$ws.Cells.Item(27, 5).Value = "  -0.22%  "
$ws.Cells.Item(28, 5).Value = "  +4.28%  "
$ws.Cells.Item(29, 4).Value = "'15.70"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.94%  "
$ws.Cells.Item(30, 4).Value = "'0.0515"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.62%  "
$ws.Cells.Item(31, 4).Value = "'1.17"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.96%  "
$ws.Cells.Item(32, 5).Value = "  +1.48%  "
$ws.Cells.Item(33, 5).Value = "  +2.72%  "
$ws.Cells.Item(34, 4).Value = "1.278.45"
$ws.Cells.Item(34, 5).Value = "  +4.95%  "
$ws.Cells.Item(35, 4).Value = "'1.54"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.23%  "
$ws.Cells.Item(36, 5).Value = "  +1.42%  "
$ws.Cells.Item(37, 5).Value = "  +2.70%  "
$ws.Cells.Item(38, 5).Value = "  +6.55%  "
$ws.Cells.Item(39, 4).Value = "'0.829"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.92%  "
$ws.Cells.Item(40, 4).Value = "'1.00"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.08%  "
$ws.Cells.Item(41, 5).Value = "  +2.57%  "
$ws.Cells.Item(42, 5).Value = "  -0.50%  "
$ws.Cells.Item(43, 5).Value = "  +1.44%  "
$ws.Cells.Item(44, 4).Value = "1.782.42"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "
$ws.Cells.Item(45, 4).Value = "'91.74"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.63%  "
$ws.Cells.Item(46, 4).Value = "'59.65"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +8.59%  "
$ws.Cells.Item(47, 4).Value = "'1.59"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.41%  "
$ws.Cells.Item(48, 4).Value = "'0.0515"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.96%  "
$ws.Cells.Item(49, 4).Value = "'7.78"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.85%  "
$ws.Cells.Item(50, 4).Value = "'0.0964"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.10%  "
$ws.Cells.Item(51, 5).Value = "  -0.52%  "
